$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 10997
$ws.Range("J51").Value = 10996
$ws.Range("L51").Value = 10996
$ws.Range("N51").Value = -11964
$ws.Range("H53").Value = 2744.1667
$ws.Range("I53").Value = 175
$ws.Range("K53").Value = 175
$ws.Range("M53").Value = 462
$ws.Range("H74").Value = 2485.9
$ws.Range("I74").Value = 2485.9
$ws.Range("K74").Value = 2485.9
$ws.Range("M74").Value = -1549.9
$ws.Range("H77").Value = 2485.9
$ws.Range("I77").Value = 2485.9
$ws.Range("K77").Value = 12429.5
$ws.Range("M77").Value = -7749.5
$ws.Range("H98").Value = 719.8
$ws.Range("I98").Value = 640.7059
$ws.Range("K98").Value = 640.7059
$ws.Range("M98").Value = 857.2941
$ws.Range("H107").Value = 951.3913
$ws.Range("I107").Value = 1135.6
$ws.Range("J107").Value = 606
$ws.Range("K107").Value = 1135.6
$ws.Range("L107").Value = 606
$ws.Range("M107").Value = 784.4000000000001
$ws.Range("N107").Value = -4446
$ws.Range("H122").Value = 719.8
$ws.Range("I122").Value = 640.7059
$ws.Range("K122").Value = 1922.1177
$ws.Range("M122").Value = 527.8822999999998
$ws.Range("H127").Value = 1186.8572
$ws.Range("I127").Value = 564.6667
$ws.Range("J127").Value = 1653.5
$ws.Range("K127").Value = 1694.0001
$ws.Range("L127").Value = 4960.5
$ws.Range("M127").Value = 3265.9999
$ws.Range("N127").Value = -14880.5
$ws.Range("H129").Value = 787.05554
$ws.Range("J129").Value = 899
$ws.Range("L129").Value = 2697
$ws.Range("N129").Value = -12697
$ws.Range("H132").Value = 2155.3057
$ws.Range("I132").Value = 2212.6287
$ws.Range("J132").Value = 149
$ws.Range("K132").Value = 6637.886100000001
$ws.Range("L132").Value = 447
$ws.Range("M132").Value = -4107.886100000001
$ws.Range("N132").Value = -5507
$ws.Range("H137").Value = 80060.84
$ws.Range("I137").Value = 4460
$ws.Range("K137").Value = 13380
$ws.Range("M137").Value = -10830
$ws.Range("H138").Value = 2269.1135
$ws.Range("J138").Value = 3498.3157
$ws.Range("L138").Value = 10494.9471
$ws.Range("N138").Value = -20774.9471

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").ClearContents()
$ws.Range("N27").Value = 0
$ws.Range("H32").Value = 26753.488
$ws.Range("I32").Value = 28263
$ws.Range("J32").Value = 7633
$ws.Range("K32").Value = 28263
$ws.Range("L32").Value = 7633
$ws.Range("M32").Value = -27976
$ws.Range("N32").Value = -8207
$ws.Range("H61").Value = 4877.778
$ws.Range("I61").Value = 2450
$ws.Range("J61").Value = 5571.4287
$ws.Range("K61").Value = 2450
$ws.Range("L61").Value = 5571.4287
$ws.Range("M61").Value = -2238
$ws.Range("N61").Value = -5995.4287
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H102").Value = 2585
$ws.Range("I102").Value = 1140
$ws.Range("J102").Value = 5475
$ws.Range("K102").Value = 1140
$ws.Range("L102").Value = 5475
$ws.Range("M102").Value = 482
$ws.Range("N102").Value = -8719
$ws.Range("H114").Value = 37039.8
$ws.Range("J114").Value = 37039.8
$ws.Range("L114").Value = 37039.8
$ws.Range("N114").Value = -45717.8
$ws.Range("H136").Value = 4877.778
$ws.Range("I136").Value = 2450
$ws.Range("J136").Value = 5571.4287
$ws.Range("K136").Value = 7350
$ws.Range("L136").Value = 16714.2861
$ws.Range("M136").Value = -4800
$ws.Range("N136").Value = -21814.2861

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2298.818
$ws.Range("I99").Value = 1929.7142
$ws.Range("J99").Value = 2944.75
$ws.Range("K99").Value = 1929.7142
$ws.Range("L99").Value = 2944.75
$ws.Range("M99").Value = -431.7141999999999
$ws.Range("N99").Value = -5940.75
$ws.Range("H107").Value = 1236.875
$ws.Range("I107").Value = 808.0909
$ws.Range("K107").Value = 808.0909
$ws.Range("M107").Value = 1111.9091
$ws.Range("H134").Value = 56648.367
$ws.Range("I134").Value = 59739.945
$ws.Range("K134").Value = 179219.835
$ws.Range("M134").Value = -176684.835

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13991.117
$ws.Range("I31").Value = 21131.842
$ws.Range("J31").Value = 4946.2
$ws.Range("K31").Value = 21131.842
$ws.Range("L31").Value = 4946.2
$ws.Range("M31").Value = -20836.842
$ws.Range("N31").Value = -5536.2
$ws.Range("H34").Value = 13991.117
$ws.Range("I34").Value = 21131.842
$ws.Range("J34").Value = 4946.2
$ws.Range("K34").Value = 21131.842
$ws.Range("L34").Value = 4946.2
$ws.Range("M34").Value = -20929.842
$ws.Range("N34").Value = -5350.2
$ws.Range("H50").Value = 15316.667
$ws.Range("J50").Value = 15316.667
$ws.Range("L50").Value = 15316.667
$ws.Range("N50").Value = -16566.667
$ws.Range("H58").Value = 37333.5
$ws.Range("I58").Value = 2285.5715
$ws.Range("J58").Value = 72381.42999999999
$ws.Range("K58").Value = 2285.5715
$ws.Range("L58").Value = 72381.42999999999
$ws.Range("M58").Value = -2082.5715
$ws.Range("N58").Value = -72787.42999999999
$ws.Range("H132").Value = 20778.074
$ws.Range("I132").Value = 26053.45
$ws.Range("K132").Value = 78160.35000000001
$ws.Range("M132").Value = -75630.35000000001
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").ClearContents()
$ws.Range("N133").Value = 0
$ws.Range("H134").Value = 1208
$ws.Range("I134").Value = 999.5
$ws.Range("J134").Value = 1625
$ws.Range("K134").Value = 2998.5
$ws.Range("L134").Value = 4875
$ws.Range("M134").Value = -463.5
$ws.Range("N134").Value = -9945
$ws.Range("H135").Value = 50395.6
$ws.Range("J135").Value = 50395.6
$ws.Range("L135").Value = 50395.6
$ws.Range("N135").Value = -60535.6
$ws.Range("H136").Value = 37333.5
$ws.Range("I136").Value = 2285.5715
$ws.Range("J136").Value = 72381.42999999999
$ws.Range("K136").Value = 6856.7145
$ws.Range("L136").Value = 217144.29
$ws.Range("M136").Value = -4306.7145
$ws.Range("N136").Value = -222244.29

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 10000
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 10000
$ws.Range("K82").Value = 0
$ws.Range("L82").ClearContents()
$ws.Range("M82").Value = 30000
$ws.Range("N82").Value = -30812
$ws.Range("H85").Value = 10000
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 10000
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -32808
$ws.Range("H131").Value = 766.12
$ws.Range("J131").Value = 771.2449
$ws.Range("L131").Value = 2313.7347
$ws.Range("N131").Value = -12393.7347

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2113.138
$ws.Range("I102").Value = 2277.5833
$ws.Range("J102").Value = 1323.8
$ws.Range("K102").Value = 2277.5833
$ws.Range("L102").Value = 1323.8
$ws.Range("M102").Value = -655.5832999999998
$ws.Range("N102").Value = -4567.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5259.9
$ws.Range("I7").Value = 5185.7144
$ws.Range("K7").Value = 5185.7144
$ws.Range("M7").Value = -5073.7144
$ws.Range("H122").Value = 1092880
$ws.Range("I122").Value = 1785512.9
$ws.Range("K122").Value = 5356538.699999999
$ws.Range("M122").Value = -5354088.699999999
$ws.Range("H124").Value = 35429
$ws.Range("J124").Value = 35429
$ws.Range("L124").Value = 35429
$ws.Range("N124").Value = -45249
$ws.Range("H126").Value = 5259.9
$ws.Range("I126").Value = 5185.7144
$ws.Range("K126").Value = 15557.1432
$ws.Range("M126").Value = -13087.1432
$ws.Range("H136").Value = 41507.54
$ws.Range("I136").Value = 52159.8
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 156479.4
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -23100
